$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1704
    $ws.Range("F3").Value = 7866
}

$wb.Save()
